$wb = $excel.ActiveWorkbook

# --- Hoja1: add rows 15-18 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$ws1.Range("A15").Value = "A004"
$ws1.Range("B15").Value = 70
$ws1.Range("C15").Value = 1000375

$ws1.Range("A16").Value = "A004"
$ws1.Range("B16").Value = 70
$ws1.Range("C16").Value = 1000376

$ws1.Range("A17").Value = "A004"
$ws1.Range("B17").Value = 70
$ws1.Range("C17").Value = 1000378

$ws1.Range("A18").Value = "A004"
$ws1.Range("B18").Value = 70
$ws1.Range("C18").Value = 1000379

# --- Hoja2: add rows 13-17 ---
$ws2 = $wb.Worksheets.Item("Hoja2")

$ws2.Range("A13").Value = 100292
$ws2.Range("B13").Value = 2024
$ws2.Range("C13").Value = "A004"
$ws2.Range("D13").Value = 1000392

$ws2.Range("A14").Value = 100293
$ws2.Range("B14").Value = 2024
$ws2.Range("C14").Value = "A004"
$ws2.Range("D14").Value = 1000394

$ws2.Range("A15").Value = 100294
$ws2.Range("B15").Value = 2024
$ws2.Range("C15").Value = "A004"
$ws2.Range("D15").Value = 1000395

$ws2.Range("A16").Value = 100295
$ws2.Range("B16").Value = 2024
$ws2.Range("C16").Value = "A004"
$ws2.Range("D16").Value = 1000396

$ws2.Range("A17").Value = 100296
$ws2.Range("B17").Value = 2024
$ws2.Range("C17").Value = "A004"
$ws2.Range("D17").Value = 1000397

$wb.Save()
